$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 646 (the post about "あす日が昇るだろう..."),
# causing all subsequent rows to shift up by one.
$ws.Rows.Item(646).Delete()
